{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// The receipt document has three paragraphs whose numbers need updating:\n//   \"\u043d\u0430 \u0441\u0443\u043c\u043c\u0443 123 USD,\"                                        -> \"\u043d\u0430 \u0441\u0443\u043c\u043c\u0443 60 USD,\"\n//   \"\u043e\u0442\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u043d\u0443\u044e \u0432\u0438\u0434\u0430\u043c \u043f\u043e\u0434 \u0441\u043b\u0435\u0434\u0443\u044e\u0449\u0438\u043c\u0438 \u0438\u043d\u0434\u0435\u043a\u0441\u0430\u043c\u0438: 3, 2,\"       -> \"\u043e\u0442\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u043d\u0443\u044e \u0432\u0438\u0434\u0430\u043c \u043f\u043e\u0434 \u0441\u043b\u0435\u0434\u0443\u044e\u0449\u0438\u043c\u0438 \u0438\u043d\u0434\u0435\u043a\u0441\u0430\u043c\u0438: 12, 7, 11,\"\n//   \"\u043d\u0430\u0445\u043e\u0434\u044f\u0449\u0438\u043c\u0441\u044f \u0432 \u0437\u0430\u043f\u043e\u0432\u0435\u0434\u043d\u0438\u043a\u0435 \u0441 \u0438\u043d\u0434\u0435\u043a\u0441\u043e\u043c 6\"                   -> \"\u043d\u0430\u0445\u043e\u0434\u044f\u0449\u0438\u043c\u0441\u044f \u0432 \u0437\u0430\u043f\u043e\u0432\u0435\u0434\u043d\u0438\u043a\u0435 \u0441 \u0438\u043d\u0434\u0435\u043a\u0441\u043e\u043c 1\"\n\nconst replacements = [\n  { find: \"\u043d\u0430 \u0441\u0443\u043c\u043c\u0443 123 USD,\", replace: \"\u043d\u0430 \u0441\u0443\u043c\u043c\u0443 60 USD,\" },\n  {\n    find: \"\u043e\u0442\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u043d\u0443\u044e \u0432\u0438\u0434\u0430\u043c \u043f\u043e\u0434 \u0441\u043b\u0435\u0434\u0443\u044e\u0449\u0438\u043c\u0438 \u0438\u043d\u0434\u0435\u043a\u0441\u0430\u043c\u0438: 3, 2,\",\n    replace: \"\u043e\u0442\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u043d\u0443\u044e \u0432\u0438\u0434\u0430\u043c \u043f\u043e\u0434 \u0441\u043b\u0435\u0434\u0443\u044e\u0449\u0438\u043c\u0438 \u0438\u043d\u0434\u0435\u043a\u0441\u0430\u043c\u0438: 12, 7, 11,\",\n  },\n  { find: \"\u043d\u0430\u0445\u043e\u0434\u044f\u0449\u0438\u043c\u0441\u044f \u0432 \u0437\u0430\u043f\u043e\u0432\u0435\u0434\u043d\u0438\u043a\u0435 \u0441 \u0438\u043d\u0434\u0435\u043a\u0441\u043e\u043c 6\", replace: \"\u043d\u0430\u0445\u043e\u0434\u044f\u0449\u0438\u043c\u0441\u044f \u0432 \u0437\u0430\u043f\u043e\u0432\u0435\u0434\u043d\u0438\u043a\u0435 \u0441 \u0438\u043d\u0434\u0435\u043a\u0441\u043e\u043c 1\" },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${find}\"`);\n  }\n\n  results.items[0].insertText(replace, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d.\n#\n# The receipt document has three paragraphs whose numbers need updating:\n#   \"\u043d\u0430 \u0441\u0443\u043c\u043c\u0443 123 USD,\"                                        -> \"\u043d\u0430 \u0441\u0443\u043c\u043c\u0443 60 USD,\"\n#   \"\u043e\u0442\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u043d\u0443\u044e \u0432\u0438\u0434\u0430\u043c \u043f\u043e\u0434 \u0441\u043b\u0435\u0434\u0443\u044e\u0449\u0438\u043c\u0438 \u0438\u043d\u0434\u0435\u043a\u0441\u0430\u043c\u0438: 3, 2,\"       -> \"\u043e\u0442\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u043d\u0443\u044e \u0432\u0438\u0434\u0430\u043c \u043f\u043e\u0434 \u0441\u043b\u0435\u0434\u0443\u044e\u0449\u0438\u043c\u0438 \u0438\u043d\u0434\u0435\u043a\u0441\u0430\u043c\u0438: 12, 7, 11,\"\n#   \"\u043d\u0430\u0445\u043e\u0434\u044f\u0449\u0438\u043c\u0441\u044f \u0432 \u0437\u0430\u043f\u043e\u0432\u0435\u0434\u043d\u0438\u043a\u0435 \u0441 \u0438\u043d\u0434\u0435\u043a\u0441\u043e\u043c 6\"                   -> \"\u043d\u0430\u0445\u043e\u0434\u044f\u0449\u0438\u043c\u0441\u044f \u0432 \u0437\u0430\u043f\u043e\u0432\u0435\u0434\u043d\u0438\u043a\u0435 \u0441 \u0438\u043d\u0434\u0435\u043a\u0441\u043e\u043c 1\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $found = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"Could not find text to replace: $findText\"\n    }\n}\n\nReplace-Text \"\u043d\u0430 \u0441\u0443\u043c\u043c\u0443 123 USD,\" \"\u043d\u0430 \u0441\u0443\u043c\u043c\u0443 60 USD,\"\nReplace-Text \"\u043e\u0442\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u043d\u0443\u044e \u0432\u0438\u0434\u0430\u043c \u043f\u043e\u0434 \u0441\u043b\u0435\u0434\u0443\u044e\u0449\u0438\u043c\u0438 \u0438\u043d\u0434\u0435\u043a\u0441\u0430\u043c\u0438: 3, 2,\" \"\u043e\u0442\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u043d\u0443\u044e \u0432\u0438\u0434\u0430\u043c \u043f\u043e\u0434 \u0441\u043b\u0435\u0434\u0443\u044e\u0449\u0438\u043c\u0438 \u0438\u043d\u0434\u0435\u043a\u0441\u0430\u043c\u0438: 12, 7, 11,\"\nReplace-Text \"\u043d\u0430\u0445\u043e\u0434\u044f\u0449\u0438\u043c\u0441\u044f \u0432 \u0437\u0430\u043f\u043e\u0432\u0435\u0434\u043d\u0438\u043a\u0435 \u0441 \u0438\u043d\u0434\u0435\u043a\u0441\u043e\u043c 6\" \"\u043d\u0430\u0445\u043e\u0434\u044f\u0449\u0438\u043c\u0441\u044f \u0432 \u0437\u0430\u043f\u043e\u0432\u0435\u0434\u043d\u0438\u043a\u0435 \u0441 \u0438\u043d\u0434\u0435\u043a\u0441\u043e\u043c 1\"\n"}
